# Update TPM-derived NATMI ligand-receptor metrics (rows 2-10) to match
# the refreshed TPM input. Only numeric <v> cells change; no structural,
# string-table, or style edits are required.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.022792
$ws.Range("H2").Value = 0.068376
$ws.Range("I2").Value = 0.001916327914826657
$ws.Range("J2").Value = 0.001916327914826657
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 0.002519701184
$ws.Range("R2").Value = 0.022677310656
$ws.Range("S2").Value = 0.00002158263474194613
$ws.Range("T2").Value = 0.00002158263474194613
# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.022792
$ws.Range("H3").Value = 0.068376
$ws.Range("I3").Value = 0.001916327914826657
$ws.Range("J3").Value = 0.001916327914826657
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("Q3").Value = 0.2054031146346667
$ws.Range("R3").Value = 1.848628031712
$ws.Range("S3").Value = 0.001759391322339475
$ws.Range("T3").Value = 0.001759391322339475
# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.022792
$ws.Range("H4").Value = 0.068376
$ws.Range("I4").Value = 0.001916327914826657
$ws.Range("J4").Value = 0.001916327914826657
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 0.015802126648
$ws.Range("R4").Value = 0.142219139832
$ws.Range("S4").Value = 0.0001353539577452361
$ws.Range("T4").Value = 0.0001353539577452362
# Row 5
$ws.Range("I5").Value = 0.3701235913233977
$ws.Range("J5").Value = 0.3701235913233977
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 0.4866603695893333
$ws.Range("R5").Value = 4.379943326304
$ws.Range("S5").Value = 0.004168515325120032
$ws.Range("T5").Value = 0.004168515325120032
# Row 6
$ws.Range("I6").Value = 0.3701235913233977
$ws.Range("J6").Value = 0.3701235913233977
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("S6").Value = 0.339812528810557
$ws.Range("T6").Value = 0.3398125288105571
# Row 7
$ws.Range("I7").Value = 0.3701235913233977
$ws.Range("J7").Value = 0.3701235913233977
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 3.052055872198666
$ws.Range("R7").Value = 27.468502849788
$ws.Range("S7").Value = 0.02614254718772066
$ws.Range("T7").Value = 0.02614254718772067
# Row 8
$ws.Range("G8").Value = 7.468693666666667
$ws.Range("H8").Value = 22.406081
$ws.Range("I8").Value = 0.6279600807617757
$ws.Range("J8").Value = 0.6279600807617757
$ws.Range("M8").Value = 0.110552
$ws.Range("N8").Value = 0.331656
$ws.Range("O8").Value = 0.01126249561724847
$ws.Range("P8").Value = 0.01126249561724847
$ws.Range("Q8").Value = 0.8256790222373334
$ws.Range("R8").Value = 7.431111200136001
$ws.Range("S8").Value = 0.007072397657386495
$ws.Range("T8").Value = 0.007072397657386496
# Row 9
$ws.Range("G9").Value = 7.468693666666667
$ws.Range("H9").Value = 22.406081
$ws.Range("I9").Value = 0.6279600807617757
$ws.Range("J9").Value = 0.6279600807617757
$ws.Range("O9").Value = 0.9181055646724333
$ws.Range("P9").Value = 0.9181055646724334
$ws.Range("Q9").Value = 67.30839511168578
$ws.Range("R9").Value = 605.775556005172
$ws.Range("S9").Value = 0.5765336445395369
$ws.Range("T9").Value = 0.5765336445395369
# Row 10
$ws.Range("G10").Value = 7.468693666666667
$ws.Range("H10").Value = 22.406081
$ws.Range("I10").Value = 0.6279600807617757
$ws.Range("J10").Value = 0.6279600807617757
$ws.Range("M10").Value = 0.6933189999999999
$ws.Range("N10").Value = 2.079957
$ws.Range("O10").Value = 0.07063193971031816
$ws.Range("P10").Value = 0.07063193971031817
$ws.Range("Q10").Value = 5.178187224279666
$ws.Range("R10").Value = 46.603685018517
$ws.Range("S10").Value = 0.04435403856485226
$ws.Range("T10").Value = 0.04435403856485227
